$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add UrbanRatioModelled value (new cell G2)
$ws.Range("G2").Value = 3332.460923935724

# Fix PopEndYear (H2)
$ws.Range("H2").Value = 0.1598809861329305

# Fix MinNightLights (V2)
$ws.Range("V2").Value = 0.7123034778733883

# Fix DistToTrans (W2)
$ws.Range("W2").Value = 60

# MaxGridDist (X2) stays the same (1), no change needed

# Fix MaxRoadDist (Y2)
$ws.Range("Y2").Value = 5

# Fix PopCutOffRoundOne (Z2)
$ws.Range("Z2").Value = 0.5
